# Update retention metrics for Sheet1 (metricas_retencao_anual)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 36: cohort_year=2019, period_index=1 -> num_customers 146 -> 147, retention_rate recalculated
$ws.Range("C36").Value = 147
$ws.Range("E36").Value = 147 / 1930

# Row 37: cohort_year=2020, period_index=0 -> num_customers/cohort_size 974 -> 978
$ws.Range("C37").Value = 978
$ws.Range("D37").Value = 978
